# İş Takip Güncellemesi - 03.08.2025 17:53:04
# Add a hidden "__id" identifier column to the two data sheets:
#   - "İş Takip Listesi" (sheet1): column M, ids 0..120 for rows 2..122
#   - "Güncelleme"        (sheet2): column Q, ids 121..148 for rows 2..29
# The id sequence is continuous across both sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("İş Takip Listesi")
$ws1.Cells.Item(1, 13).Value = "__id"
$id = 0
for ($row = 2; $row -le 122; $row++) {
    $ws1.Cells.Item($row, 13).Value = $id
    $id = $id + 1
}

$ws2 = $wb.Worksheets.Item("Güncelleme")
$ws2.Cells.Item(1, 17).Value = "__id"
for ($row = 2; $row -le 29; $row++) {
    $ws2.Cells.Item($row, 17).Value = $id
    $id = $id + 1
}
